$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.011056780815125
$ws.Range("B1").Value = 1.259087800979614
$ws.Range("C1").Value = 5.613296508789062
$ws.Range("D1").Value = 1.653684496879578
$ws.Range("E1").Value = 1.011772513389587
